$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.787217997188553
$ws.Range("D2").Value = 4.060893136428637
$ws.Range("E2").Value = 13.20799405628988
$ws.Range("F2").Value = 20.15661207283617
$ws.Range("G2").Value = 21.77006372168183
$ws.Range("H2").Value = 12.43018211106743
$ws.Range("K2").Value = 10.64280726498
$ws.Range("M2").Value = 12.85440798663801
$ws.Range("N2").Value = 17.86358411131111
$ws.Range("O2").Value = 17.99543016889661
$ws.Range("B3").Value = 7.714374659171505
$ws.Range("D3").Value = 3.982170352681993
$ws.Range("E3").Value = 13.06334835640444
$ws.Range("F3").Value = 20.13270166463773
$ws.Range("G3").Value = 21.70534967497087
$ws.Range("H3").Value = 12.4643365330203
$ws.Range("K3").Value = 10.31498380857928
$ws.Range("M3").Value = 12.6090107184405
$ws.Range("N3").Value = 17.91829416268558
$ws.Range("O3").Value = 18.03233334819969
$ws.Range("B4").Value = 7.671135964444422
$ws.Range("D4").Value = 3.932372254637404
$ws.Range("E4").Value = 12.9786258501203
$ws.Range("F4").Value = 20.12391451323946
$ws.Range("G4").Value = 21.67382517052684
$ws.Range("H4").Value = 12.48740343930805
$ws.Range("K4").Value = 10.10654140856615
$ws.Range("M4").Value = 12.45871070603053
$ws.Range("N4").Value = 17.95363580648078
$ws.Range("O4").Value = 18.05930115476171
$ws.Range("B5").Value = 7.653909753706696
$ws.Range("D5").Value = 3.911725813124971
$ws.Range("E5").Value = 12.9451706429843
$ws.Range("F5").Value = 20.12181902250849
$ws.Range("G5").Value = 21.66305259456006
$ws.Range("H5").Value = 12.49733009322655
$ws.Range("K5").Value = 10.01988045988559
$ws.Range("M5").Value = 12.39764256042835
$ws.Range("N5").Value = 17.96847883282249
$ws.Range("O5").Value = 18.07137192667942
$ws.Range("B6").Value = 7.651073712385926
$ws.Range("D6").Value = 3.908276580801455
$ws.Range("E6").Value = 12.93968120511576
$ws.Range("F6").Value = 20.1215608470428
$ws.Range("G6").Value = 21.6613892895046
$ws.Range("H6").Value = 12.49901020761226
$ws.Range("K6").Value = 10.00538916448124
$ws.Range("M6").Value = 12.38751567023416
$ws.Range("N6").Value = 17.97097017520829
$ws.Range("O6").Value = 18.07344149377884
$ws.Range("B7").Value = 7.67090202508406
$ws.Range("D7").Value = 3.932095220610515
$ws.Range("E7").Value = 12.97817027692089
$ws.Range("F7").Value = 20.12388023534945
$ws.Range("G7").Value = 21.67367148050949
$ws.Range("H7").Value = 12.48753518151254
$ws.Range("K7").Value = 10.10537951641961
$ws.Range("M7").Value = 12.45788627495738
$ws.Range("N7").Value = 17.95383419758197
$ws.Range("O7").Value = 18.05945957117327
$ws.Range("B8").Value = 7.761805050149243
$ws.Range("D8").Value = 4.034061313761114
$ws.Range("E8").Value = 13.15729680622191
$ws.Range("F8").Value = 20.14714685568606
$ws.Range("G8").Value = 21.74605287961489
$ws.Range("H8").Value = 12.4415233423642
$ws.Range("K8").Value = 10.53130542276325
$ws.Range("M8").Value = 12.76976841411508
$ws.Range("N8").Value = 17.88208570464177
$ws.Range("O8").Value = 18.00725867938297
$ws.Range("B9").Value = 7.95096134100243
$ws.Range("D9").Value = 4.221826639691846
$ws.Range("E9").Value = 13.53908268893657
$ws.Range("F9").Value = 20.23933218766745
$ws.Range("G9").Value = 21.95257575815238
$ws.Range("H9").Value = 12.36793942349507
$ws.Range("K9").Value = 11.30658346900733
$ws.Range("M9").Value = 13.38064248040476
$ws.Range("N9").Value = 17.75521754504648
$ws.Range("O9").Value = 17.93917901459272
$ws.Range("B10").Value = 8.095326893460918
$ws.Range("D10").Value = 4.351613293197904
$ws.Range("E10").Value = 13.83545052939753
$ws.Range("F10").Value = 20.33507721625972
$ws.Range("G10").Value = 22.14271835082421
$ws.Range("H10").Value = 12.32404384199839
$ws.Range("K10").Value = 11.83588863876779
$ws.Range("M10").Value = 13.8242161237803
$ws.Range("N10").Value = 17.67036762634267
$ws.Range("O10").Value = 17.91017564407907
$ws.Range("B11").Value = 8.161901133953082
$ws.Range("D11").Value = 4.40874432500953
$ws.Range("E11").Value = 13.97311196271167
$ws.Range("F11").Value = 20.38460964311158
$ws.Range("G11").Value = 22.23729246676171
$ws.Range("H11").Value = 12.3062857763092
$ws.Range("K11").Value = 12.0672784324167
$ws.Range("M11").Value = 14.02391787625523
$ws.Range("N11").Value = 17.63356749561506
$ws.Range("O11").Value = 17.90156180530126
$ws.Range("B12").Value = 8.187217759333207
$ws.Range("D12").Value = 4.430093115426159
$ws.Range("E12").Value = 14.02559746443209
$ws.Range("F12").Value = 20.40421450016228
$ws.Range("G12").Value = 22.27424003910117
$ws.Range("H12").Value = 12.29987934806285
$ws.Range("K12").Value = 12.15350219556858
$ws.Range("M12").Value = 14.09916573846776
$ws.Range("N12").Value = 17.61988977043013
$ws.Range("O12").Value = 17.89895935672618
$ws.Range("B13").Value = 8.181761007750547
$ws.Range("D13").Value = 4.425508133945309
$ws.Range("E13").Value = 14.01427877295551
$ws.Range("F13").Value = 20.39995472682632
$ws.Range("G13").Value = 22.26623272406581
$ws.Range("H13").Value = 12.30124493156282
$ws.Range("K13").Value = 12.13499526599669
$ws.Range("M13").Value = 14.08297759418605
$ws.Range("N13").Value = 17.62282406742208
$ws.Range("O13").Value = 17.89949050320486
$ws.Range("B14").Value = 8.163981948248225
$ws.Range("D14").Value = 4.410506484161057
$ws.Range("E14").Value = 13.97742314124148
$ws.Range("F14").Value = 20.38620562232783
$ws.Range("G14").Value = 22.24030959494007
$ws.Range("H14").Value = 12.30575233729877
$ws.Range("K14").Value = 12.07440039540143
$ws.Range("M14").Value = 14.030116493343
$ws.Range("N14").Value = 17.63243706185154
$ws.Range("O14").Value = 17.90133448098816
$ws.Range("B15").Value = 8.153104929682831
$ws.Range("D15").Value = 4.401280051019747
$ws.Range("E15").Value = 13.95489281513539
$ws.Range("F15").Value = 20.37789396737011
$ws.Range("G15").Value = 22.22457781477013
$ws.Range("H15").Value = 12.30855469734615
$ws.Range("K15").Value = 12.03710083705679
$ws.Range("M15").Value = 13.99768655639185
$ws.Range("N15").Value = 17.63835882918245
$ws.Range("O15").Value = 17.90254986559115
$ws.Range("B16").Value = 8.09099226437427
$ws.Range("D16").Value = 4.347840271028377
$ws.Range("E16").Value = 13.82650673054911
$ws.Range("F16").Value = 20.33195948202943
$ws.Range("G16").Value = 22.1366979890109
$ws.Range("H16").Value = 12.32524887867687
$ws.Range("K16").Value = 11.82057353675828
$ws.Range("M16").Value = 13.81111670377152
$ws.Range("N16").Value = 17.67280871159539
$ws.Range("O16").Value = 17.9108308188268
$ws.Range("B17").Value = 8.053102525249603
$ws.Range("D17").Value = 4.314559607467458
$ws.Range("E17").Value = 13.74843662365704
$ws.Range("F17").Value = 20.30530336455144
$ws.Range("G17").Value = 22.08483690784952
$ws.Range("H17").Value = 12.33605656615687
$ws.Range("K17").Value = 11.68529931959532
$ws.Range("M17").Value = 13.69607388666395
$ws.Range("N17").Value = 17.69440256952928
$ws.Range("O17").Value = 17.91708466171456
$ws.Range("B18").Value = 8.03139560576215
$ws.Range("D18").Value = 4.295238313737851
$ws.Range("E18").Value = 13.70380345844987
$ws.Range("F18").Value = 20.29053488483981
$ws.Range("G18").Value = 22.05576985326861
$ws.Range("H18").Value = 12.34248088394853
$ws.Range("K18").Value = 11.60661145776712
$ws.Range("M18").Value = 13.62971258315733
$ws.Range("N18").Value = 17.70699212606759
$ws.Range("O18").Value = 17.92111276598528
$ws.Range("B19").Value = 8.024061539815314
$ws.Range("D19").Value = 4.28866603376473
$ws.Range("E19").Value = 13.68873952764592
$ws.Range("F19").Value = 20.28563162801953
$ws.Range("G19").Value = 22.04605991657665
$ws.Range("H19").Value = 12.34469176477261
$ws.Range("K19").Value = 11.57981921212431
$ws.Range("M19").Value = 13.60721322413681
$ws.Range("N19").Value = 17.71128384691889
$ws.Range("O19").Value = 17.92255061611791
$ws.Range("B20").Value = 8.057127194111255
$ws.Range("D20").Value = 4.318121012271003
$ws.Range("E20").Value = 13.75671968446406
$ws.Range("F20").Value = 20.30808272190703
$ws.Range("G20").Value = 22.09027892533739
$ws.Range("H20").Value = 12.33488453735966
$ws.Range("K20").Value = 11.69979110809977
$ws.Range("M20").Value = 13.70834079258761
$ws.Range("N20").Value = 17.69208634722786
$ws.Range("O20").Value = 17.91637430918687
$ws.Range("B21").Value = 8.169201385969489
$ws.Range("D21").Value = 4.414920663474791
$ws.Range("E21").Value = 13.98823928870123
$ws.Range("F21").Value = 20.39022115194433
$ws.Range("G21").Value = 22.24789328759816
$ws.Range("H21").Value = 12.30441976555901
$ws.Range("K21").Value = 12.09223688603605
$ws.Range("M21").Value = 14.04565383435875
$ws.Range("N21").Value = 17.62960650625685
$ws.Range("O21").Value = 17.90077495862216
$ws.Range("B22").Value = 8.243057873433564
$ws.Range("D22").Value = 4.476516261964663
$ws.Range("E22").Value = 14.14160201413161
$ws.Range("F22").Value = 20.44883985682452
$ws.Range("G22").Value = 22.35750123131694
$ws.Range("H22").Value = 12.28636391944075
$ws.Range("K22").Value = 12.34055289917988
$ws.Range("M22").Value = 14.26389041001222
$ws.Range("N22").Value = 17.59027388332768
$ws.Range("O22").Value = 17.89442363457563
$ws.Range("B23").Value = 8.203591153482806
$ws.Range("D23").Value = 4.443797605031003
$ws.Range("E23").Value = 14.05957912845909
$ws.Range("F23").Value = 20.41710642960971
$ws.Range("G23").Value = 22.29840731419723
$ws.Range("H23").Value = 12.29583086330886
$ws.Range("K23").Value = 12.20878386578556
$ws.Range("M23").Value = 14.14764014720974
$ws.Range("N23").Value = 17.61112935729214
$ws.Range("O23").Value = 17.89746157431969
$ws.Range("B24").Value = 8.055307401589864
$ws.Range("D24").Value = 4.316511485269761
$ws.Range("E24").Value = 13.75297412876998
$ws.Range("F24").Value = 20.30682444019621
$ws.Range("G24").Value = 22.08781625544923
$ws.Range("H24").Value = 12.33541375486012
$ws.Range("K24").Value = 11.69324222285372
$ws.Range("M24").Value = 13.70279560903291
$ws.Range("N24").Value = 17.69313296626285
$ws.Range("O24").Value = 17.91669411183094
$ws.Range("B25").Value = 7.898750654859183
$ws.Range("D25").Value = 4.172414915369404
$ws.Range("E25").Value = 13.43280849394993
$ws.Range("F25").Value = 20.20944135069506
$ws.Range("G25").Value = 21.8898791025117
$ws.Range("H25").Value = 12.38606174879678
$ws.Range("K25").Value = 11.10368203571238
$ws.Range("M25").Value = 13.21596853567002
$ws.Range("N25").Value = 17.78806553172516
$ws.Range("O25").Value = 17.95391281063233
